$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Step 1: insert the new "2022-Q1" sheet right before "总计".
#
# We rename the existing "总计" worksheet to "2022-Q1" (this keeps its
# original sheetId=6 and its formatting/sheetPr intact) and then make a
# full copy of it, placed right after, which we rename back to "总计"
# (giving it a fresh sheetId=7). This matches the sheetId renumbering
# seen in the target workbook.xml.
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Copy($null, $q1)
$total = $wb.Worksheets.Item("2022-Q1 (2)")
$total.Name = "总计"

# ------------------------------------------------------------------
# Step 2: populate "2022-Q1" (formerly "总计") with the new fund-
# holdings table: header row + 24 data rows, columns A-H.
# ------------------------------------------------------------------
$ws = $q1

# Extend the header-row format (style of D1) across E1:H1 first,
# then set the new header captions.
$ws.Range("D1").Copy()
$ws.Range("E1:H1").PasteSpecial(-4122)
$ws.Range("B1").Value2 = "基金代码"
$ws.Range("C1").Value2 = "基金名称"
$ws.Range("D1").Value2 = "基金规模"
$ws.Range("E1").Value2 = "股票总仓位"
$ws.Range("F1").Value2 = "仓位占比"
$ws.Range("G1").Value2 = "持有市值(亿元)"
$ws.Range("H1").Value2 = "仓位排名"

# Column-A (row index) cells use the bordered/centered style already
# carried by rows 2..6. Stamp that same format down through row 25.
$ws.Range("A6").Copy()
$ws.Range("A7:A25").PasteSpecial(-4122)

$ws.Range("A2").Value2 = 0
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value2 = "000979"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value2 = "景顺长城沪港深精选股票"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "16.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value2 = "82.61"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value2 = "8.92"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value2 = "1.4682"
$ws.Range("G2").Style = "Normal"
$ws.Range("H2").Value2 = 4

$ws.Range("A3").Value2 = 1
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value2 = "260112"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value2 = "景顺长城能源基建混合"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "16.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value2 = "60.89"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value2 = "7.94"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value2 = "1.3093"
$ws.Range("G3").Style = "Normal"
$ws.Range("H3").Value2 = 3

$ws.Range("A4").Value2 = 2
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value2 = "008850"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value2 = "景顺长城价值稳进三年定期开放灵活配置混合"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "17.06"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value2 = "69.71"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value2 = "7.29"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value2 = "1.2437"
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Value2 = 2

$ws.Range("A5").Value2 = 3
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value2 = "009098"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value2 = "景顺长城价值领航两年持有期混合"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "11.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value2 = "75.58"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value2 = "9.80"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value2 = "1.1437"
$ws.Range("G5").Style = "Normal"
$ws.Range("H5").Value2 = 5

$ws.Range("A6").Value2 = 4
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value2 = "008715"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value2 = "景顺长城价值驱动一年持有期灵活配置混合型证券投资基金"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "16.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value2 = "62.03"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value2 = "5.43"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value2 = "0.9139"
$ws.Range("G6").Style = "Normal"
$ws.Range("H6").Value2 = 4

$ws.Range("A7").Value2 = 5
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value2 = "010264"
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value2 = "鹏华成长智选混合A"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "42.07"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value2 = "72.05"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value2 = "1.63"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value2 = "0.6857"
$ws.Range("G7").Style = "Normal"
$ws.Range("H7").Value2 = 9

$ws.Range("A8").Value2 = 6
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value2 = "008060"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value2 = "景顺长城价值边际灵活配置混合"
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "4.93"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value2 = "80.78"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value2 = "8.82"
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value2 = "0.4348"
$ws.Range("G8").Style = "Normal"
$ws.Range("H8").Value2 = 4

$ws.Range("A9").Value2 = 7
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value2 = "012366"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value2 = "上投摩根安荣回报混合型证券投资基金A"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "41.47"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value2 = "21.90"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value2 = "0.95"
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value2 = "0.3940"
$ws.Range("G9").Style = "Normal"
$ws.Range("H9").Value2 = 9

$ws.Range("A10").Value2 = 8
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value2 = "159611"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value2 = "广发中证全指电力ETF"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "13.38"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value2 = "99.14"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value2 = "2.19"
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value2 = "0.2930"
$ws.Range("G10").Style = "Normal"
$ws.Range("H10").Value2 = 10

$ws.Range("A11").Value2 = 9
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value2 = "004738"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value2 = "上投摩根安隆回报混合A"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "23.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value2 = "21.18"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value2 = "1.08"
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value2 = "0.2488"
$ws.Range("G11").Style = "Normal"
$ws.Range("H11").Value2 = 7

$ws.Range("A12").Value2 = 10
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value2 = "012367"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value2 = "上投摩根安荣回报混合型证券投资基金C"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "23.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value2 = "21.90"
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value2 = "0.95"
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value2 = "0.2274"
$ws.Range("G12").Style = "Normal"
$ws.Range("H12").Value2 = 9

$ws.Range("A13").Value2 = 11
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value2 = "004823"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value2 = "上投摩根安裕回报混合A"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "11.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value2 = "30.12"
$ws.Range("E13").Style = "Normal"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value2 = "1.33"
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value2 = "0.1502"
$ws.Range("G13").Style = "Normal"
$ws.Range("H13").Value2 = 10

$ws.Range("A14").Value2 = 12
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value2 = "673110"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value2 = "西部利得新润灵活配置混合"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "5.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value2 = "81.49"
$ws.Range("E14").Style = "Normal"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value2 = "2.60"
$ws.Range("F14").Style = "Normal"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value2 = "0.1326"
$ws.Range("G14").Style = "Normal"
$ws.Range("H14").Value2 = 9

$ws.Range("A15").Value2 = 13
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value2 = "004824"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value2 = "上投摩根安裕回报混合C"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "7.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value2 = "30.12"
$ws.Range("E15").Style = "Normal"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value2 = "1.33"
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value2 = "0.0992"
$ws.Range("G15").Style = "Normal"
$ws.Range("H15").Value2 = 10

$ws.Range("A16").Value2 = 14
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value2 = "007146"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value2 = "鹏华研究智选混合"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "5.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value2 = "83.09"
$ws.Range("E16").Style = "Normal"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value2 = "1.77"
$ws.Range("F16").Style = "Normal"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value2 = "0.0961"
$ws.Range("G16").Style = "Normal"
$ws.Range("H16").Value2 = 10

$ws.Range("A17").Value2 = 15
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value2 = "011349"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value2 = "淳厚现代服务业股票A"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "3.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value2 = "81.51"
$ws.Range("E17").Style = "Normal"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value2 = "2.40"
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value2 = "0.0859"
$ws.Range("G17").Style = "Normal"
$ws.Range("H17").Value2 = 9

$ws.Range("A18").Value2 = 16
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value2 = "004739"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value2 = "上投摩根安隆回报混合C"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "7.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value2 = "21.18"
$ws.Range("E18").Style = "Normal"
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value2 = "1.08"
$ws.Range("F18").Style = "Normal"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value2 = "0.0791"
$ws.Range("G18").Style = "Normal"
$ws.Range("H18").Value2 = 7

$ws.Range("A19").Value2 = 17
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value2 = "010265"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value2 = "鹏华成长智选混合C"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "3.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value2 = "72.05"
$ws.Range("E19").Style = "Normal"
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value2 = "1.63"
$ws.Range("F19").Style = "Normal"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value2 = "0.0553"
$ws.Range("G19").Style = "Normal"
$ws.Range("H19").Value2 = 9

$ws.Range("A20").Value2 = 18
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value2 = "512390"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value2 = "平安MSCI中国A股低波动ETF"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "2.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value2 = "97.88"
$ws.Range("E20").Style = "Normal"
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value2 = "1.76"
$ws.Range("F20").Style = "Normal"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value2 = "0.0488"
$ws.Range("G20").Style = "Normal"
$ws.Range("H20").Value2 = 8

$ws.Range("A21").Value2 = 19
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value2 = "006700"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value2 = "红土创新稳健混合A"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "0.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value2 = "27.06"
$ws.Range("E21").Style = "Normal"
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value2 = "5.23"
$ws.Range("F21").Style = "Normal"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value2 = "0.0387"
$ws.Range("G21").Style = "Normal"
$ws.Range("H21").Value2 = 2

$ws.Range("A22").Value2 = 20
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value2 = "006701"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value2 = "红土创新稳健混合C"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "0.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value2 = "27.06"
$ws.Range("E22").Style = "Normal"
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value2 = "5.23"
$ws.Range("F22").Style = "Normal"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value2 = "0.0183"
$ws.Range("G22").Style = "Normal"
$ws.Range("H22").Value2 = 2

$ws.Range("A23").Value2 = 21
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value2 = "011350"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value2 = "淳厚现代服务业股票C"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "0.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value2 = "81.51"
$ws.Range("E23").Style = "Normal"
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value2 = "2.40"
$ws.Range("F23").Style = "Normal"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value2 = "0.0151"
$ws.Range("G23").Style = "Normal"
$ws.Range("H23").Value2 = 9

$ws.Range("A24").Value2 = 22
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value2 = "159962"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value2 = "华夏中证四川国企改革ETF"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "0.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value2 = "95.82"
$ws.Range("E24").Style = "Normal"
$ws.Range("F24").NumberFormat = "@"
$ws.Range("F24").Value2 = "2.88"
$ws.Range("F24").Style = "Normal"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value2 = "0.0141"
$ws.Range("G24").Style = "Normal"
$ws.Range("H24").Value2 = 10

$ws.Range("A25").Value2 = 23
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value2 = "009188"
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value2 = "鹏华股息精选混合"
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "0.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value2 = "89.69"
$ws.Range("E25").Style = "Normal"
$ws.Range("F25").NumberFormat = "@"
$ws.Range("F25").Value2 = "2.05"
$ws.Range("F25").Style = "Normal"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value2 = "0.0131"
$ws.Range("G25").Style = "Normal"
$ws.Range("H25").Value2 = 2

# ------------------------------------------------------------------
# Step 3: rebuild "总计" with a new 2022-Q1 row inserted at the top
# (index 0), pushing the previously existing rows down by one and
# renumbering their index column (A) accordingly.
# ------------------------------------------------------------------
# Stamp the A-column style down to the new last row (7) first.
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122)

# Shift existing data rows 2..6 down to rows 3..7 (bottom-up so we
# never clobber a row before it has been read).
for ($r = 6; $r -ge 2; $r--) {
    $dest = $r + 1
    $bval = $total.Range("B$r").Value2
    $cval = $total.Range("C$r").Value2
    $dval = $total.Range("D$r").Value2
    $total.Range("A$dest").Value2 = $r - 1
    $total.Range("B$dest").Value2 = $bval
    $total.Range("C$dest").Value2 = $cval
    $total.Range("D$dest").Value2 = $dval
}

# Write the new 2022-Q1 summary row at the top.
$total.Range("A2").Value2 = 0
$total.Range("B2").Value2 = "2022-Q1"
$total.Range("C2").Value2 = 24
$total.Range("D2").Value2 = 9.21
